# BOM fix: swap the op-amp part (U2) from LMV321 to AS321 so the board can
# handle the higher regulator voltage, and update its LCSC order code.
#
# Sheet layout (row 15 = U2 line item):
#   A = Comment (part description)
#   B = Designator
#   C = Footprint
#   D = LCSC Part # (optional)
#   E = Qty

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment: LMV321 Operational Amplifier -> AS321 Operational Amplifier
$ws.Range("A15").Value = "AS321 Operational Amplifier "

# LCSC Part #: C686637 -> C144156
$ws.Range("D15").Value = "C144156"

# Leave the selection on the cell that was just edited (D15), matching the
# cursor position after the change.
$ws.Range("D15").Select()

# Column D was narrowed slightly (23.5 -> 21.42 chars) alongside the edit.
$ws.Columns("D").ColumnWidth = 20.714285714285715
